$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 held a standalone numeric placeholder cell (bold, bordered, centered).
# Deleting it shifts the question-data row (previously row 2, a plain/default-
# styled cell) up into row 1, dropping the special formatting along with it.
$ws.Rows(1).Delete()

# Replace the cell text with the pretty-printed (indented, double-quoted) JSON
# rendering of the same question data.
$text = @'
questions = [
    {
        "title": "You are implementing a code for retrieving SparkSessions in a data processing pipeline. You want the code to first check if a session with the given configuration exists. If it does not exist, the code should create a new session and assign it as the global default session. Which method call fulfills these requirements?",
        "ques_type": 2,
        "options": [
            "SparkSession.builder.config(\"key\", \"val\").getOrCreate()",
            "SparkSession.builder.config(\"key\", \"val\").getOrCreate().setGlobal()",
            "SparkSession.builder.config(\"key\", \"val\").getOrCreateGlobal()",
            "SparkSession.builder.config(\"key\", \"val\").getDefaultGlobal()"
        ],
        "score": "SparkSession.builder.config(\"key\", \"val\").getOrCreate()"
    },
    {
        "title": "You are implementing a data processing pipeline using the Pandas API on Spark. You want to be able to detect non-missing values in the dataframe. Which functions can you use to accomplish this?",
        "ques_type": 15,
        "options": [
            "pyspark.pandas.nonmissing",
            "pyspark.pandas.notNA",
            "pyspark.pandas.isnotNA",
            "pyspark.pandas.notna",
            "pyspark.pandas.notnull"
        ],
        "score": [
            "pyspark.pandas.notna",
            "pyspark.pandas.notnull"
        ]
    },
    {
        "title": "You are reviewing your co-worker\u2019s code, related to machine learning with PySpark. You are currently looking at vectors, and you notice that most of the values in the vectors are 0s. You want to convert regular vectors to sparse ones and switch from the pyspark.ml.linalg.DenseVector function to pyspark.ml.linalg.SparseVector. No additional changes are required. How many arguments will you need to reassign to switch from DenseVector to SparseVector?",
        "ques_type": 2,
        "options": [
            "0",
            "1",
            "2",
            "3"
        ],
        "score": "0"
    },
    {
        "title": "You are designing a Spark application, and you want to define a cluster manager to which the application should connect when it starts. Which function from the SparkConf module should you use?",
        "ques_type": 2,
        "options": [
            "setMaster",
            "setSparkHome",
            "setExecutorEnv",
            "setAppName"
        ],
        "score": "setMaster"
    }
]
'@
$ws.Range("A1").Value = $text

# Re-fit the row height after the multi-line assignment so it reports the sheet
# default height instead of an oversized auto-computed one.
$ws.Rows(1).AutoFit()

